# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'52.533.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.90%  "

# Row 3
$ws.Range("D3").Value = "'3.114.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.13%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'394.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.67%  "

# Row 6
$ws.Range("D6").Value = "'103.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.84%  "

# Row 7
$ws.Range("D7").Value = "'0.540"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.66%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").Value = "'0.603"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.19%  "

# Row 10
$ws.Range("E10").Value = "  +3.30%  "

# Row 11
$ws.Range("E11").Value = "  +1.04%  "

# Row 12
$ws.Range("D12").Value = "'0.0862"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "

# Row 13
$ws.Range("D13").Value = "'3.614.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.03%  "

# Row 14
$ws.Range("D14").Value = "'18.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.71%  "

# Row 15
$ws.Range("D15").Value = "'7.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.29%  "

# Row 16
$ws.Range("D16").Value = "'1.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.33%  "

# Row 17
$ws.Range("D17").Value = "'3.117.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.36%  "

# Row 18
$ws.Range("D18").Value = "'10.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.79%  "

# Row 19
$ws.Range("D19").Value = "'52.450.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.59%  "

# Row 20
$ws.Range("D20").Value = "'3.23"
$ws.Range("D20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'12.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.44%  "

# Row 22
$ws.Range("D22").Value = "'0.0₃0970"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.69%  "

# Row 23
$ws.Range("D23").Value = "'70.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.95%  "

# Row 24
$ws.Range("D24").Value = "'268.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.09%  "

# Row 25
$ws.Range("D25").Value = "'3.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.55%  "

# Row 26
$ws.Range("D26").Value = "'8.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.98%  "

# Row 27
$ws.Range("D27").Value = "'27.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.03%  "

# Row 28
$ws.Range("D28").Value = "'7.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.86%  "

# Row 29
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.03%  "

# Row 30
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.167"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.71%  "

# Row 31
$ws.Range("E31").Value = "  +0.06%  "

# Row 32
$ws.Range("D32").Value = "'10.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.67%  "

# Row 33
$ws.Range("B33").Value = "VeChain"
$ws.Range("C33").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D33").Value = "'0.0493"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.86%  "

# Row 34
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "'36.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.22%  "

# Row 35
$ws.Range("E35").Value = "  +0.88%  "

# Row 36
$ws.Range("D36").Value = "'49.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.14%  "

# Row 37
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "

# Row 38
$ws.Range("D38").Value = "'3.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.94%  "

# Row 39
$ws.Range("E39").Value = "  +10.54%  "

# Row 40
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.293"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.17%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.69%  "

# Row 42
$ws.Range("D42").Value = "'130.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.52%  "

# Row 43
$ws.Range("D43").Value = "'16.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "

# Row 44
$ws.Range("D44").Value = "'1.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.25%  "

# Row 45
$ws.Range("D45").Value = "'0.116"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.45%  "

# Row 46
$ws.Range("D46").Value = "'22.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.82%  "

# Row 47
$ws.Range("D47").Value = "'2.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.85%  "

# Row 48
$ws.Range("E48").Value = "  -0.88%  "

# Row 49
$ws.Range("D49").Value = "'2.082.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.96%  "

# Row 50
$ws.Range("D50").Value = "'0.0534"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +36.52%  "

# Row 51
$ws.Range("D51").Value = "'0.914"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.61%  "

Write-Host "Updated cryptos list"